{"js": "// Map of old multiplication-problem text -> new multiplication-problem text,\n// as described by the unified OOXML diff.\nconst replacements = {\n  \"34\u00d755=\": \"50\u00d720=\",\n  \"63\u00d757=\": \"97\u00d777=\",\n  \"54\u00d773=\": \"31\u00d721=\",\n  \"45\u00d795=\": \"92\u00d769=\",\n  \"39\u00d724=\": \"52\u00d756=\",\n  \"33\u00d794=\": \"53\u00d741=\",\n  \"73\u00d728=\": \"79\u00d785=\",\n  \"47\u00d720=\": \"24\u00d792=\",\n  \"19\u00d714=\": \"48\u00d732=\",\n  \"70\u00d756=\": \"31\u00d753=\",\n  \"93\u00d729=\": \"55\u00d734=\",\n  \"97\u00d792=\": \"20\u00d731=\",\n  \"37\u00d719=\": \"13\u00d716=\",\n  \"24\u00d785=\": \"24\u00d730=\",\n  \"84\u00d782=\": \"34\u00d723=\",\n  \"16\u00d777=\": \"27\u00d782=\",\n  \"41\u00d741=\": \"42\u00d792=\",\n  \"83\u00d791=\": \"80\u00d764=\",\n  \"59\u00d758=\": \"62\u00d760=\",\n  \"47\u00d772=\": \"47\u00d797=\",\n  \"44\u00d733=\": \"14\u00d786=\",\n  \"29\u00d714=\": \"98\u00d798=\",\n  \"72\u00d793=\": \"32\u00d766=\",\n  \"94\u00d775=\": \"51\u00d743=\",\n  \"42\u00d733=\": \"90\u00d739=\",\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  // `columnCount` isn't reliably populated by `.load()` in this host, so\n  // derive the column count from a `values` snapshot instead.\n  table.load(\"rowCount,values\");\n  await context.sync();\n\n  const rowCount = table.rowCount;\n  const columnCount = table.values.length > 0 ? table.values[0].length : 0;\n\n  for (let r = 0; r < rowCount; r++) {\n    for (let c = 0; c < columnCount; c++) {\n      const cell = table.getCell(r, c);\n      cell.load(\"value\");\n      await context.sync();\n\n      const text = cell.value;\n      if (Object.prototype.hasOwnProperty.call(replacements, text)) {\n        cell.value = replacements[text];\n        await context.sync();\n      }\n    }\n  }\n}\n", "ps1": "# Map of old multiplication-problem text -> new multiplication-problem text,\n# as described by the unified OOXML diff.\n$map = @{\n  \"34\u00d755=\" = \"50\u00d720=\";\n  \"63\u00d757=\" = \"97\u00d777=\";\n  \"54\u00d773=\" = \"31\u00d721=\";\n  \"45\u00d795=\" = \"92\u00d769=\";\n  \"39\u00d724=\" = \"52\u00d756=\";\n  \"33\u00d794=\" = \"53\u00d741=\";\n  \"73\u00d728=\" = \"79\u00d785=\";\n  \"47\u00d720=\" = \"24\u00d792=\";\n  \"19\u00d714=\" = \"48\u00d732=\";\n  \"70\u00d756=\" = \"31\u00d753=\";\n  \"93\u00d729=\" = \"55\u00d734=\";\n  \"97\u00d792=\" = \"20\u00d731=\";\n  \"37\u00d719=\" = \"13\u00d716=\";\n  \"24\u00d785=\" = \"24\u00d730=\";\n  \"84\u00d782=\" = \"34\u00d723=\";\n  \"16\u00d777=\" = \"27\u00d782=\";\n  \"41\u00d741=\" = \"42\u00d792=\";\n  \"83\u00d791=\" = \"80\u00d764=\";\n  \"59\u00d758=\" = \"62\u00d760=\";\n  \"47\u00d772=\" = \"47\u00d797=\";\n  \"44\u00d733=\" = \"14\u00d786=\";\n  \"29\u00d714=\" = \"98\u00d798=\";\n  \"72\u00d793=\" = \"32\u00d766=\";\n  \"94\u00d775=\" = \"51\u00d743=\";\n  \"42\u00d733=\" = \"90\u00d739=\";\n}\n\n$d = $word.ActiveDocument\n\nforeach ($t in $d.Tables) {\n  $rows = $t.Rows.Count\n  $cols = $t.Columns.Count\n\n  for ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n      $cell = $t.Cell($r, $c)\n      $txt = $cell.Range.Text\n      # Cell ranges include the trailing cell-mark (CR + BEL); strip it\n      # before comparing against the lookup table.\n      $clean = $txt.TrimEnd([char]7, [char]13)\n      if ($map.ContainsKey($clean)) {\n        $cell.Range.Text = $map[$clean]\n      }\n    }\n  }\n}\n"}
